# Scene 57 edit: add Roxy's expression/pose tags to her dialogue lines,
# and insert new stand-alone "tag" paragraphs marking expression changes
# that occur between/around existing lines (per commit "write some new
# for stephen").

$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Output "REPLACE NOT FOUND: $oldText"
    }
}

function Insert-TagAfter($anchorText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        Write-Output "ANCHOR NOT FOUND: $anchorText"
        return
    }
    $idx = $rng.Paragraphs.Item(1).Index
    $rng.InsertParagraphAfter()
    $d.Paragraphs.Item($idx + 1).Range.Text = $newText
}

# --- 1) Update existing Roxy/dialogue lines with expression tags -----------

Replace-ExactText "Roxy: Hey, there." "Roxy (waving smiling): Hey, there."

Replace-ExactText "Roxy: I’m actually pretty relieved that you guys came…" "Roxy (neutral smiling_worried): I’m actually pretty relieved that you guys came…"

Replace-ExactText "Roxy: I haven’t played with everyone else for a while, and I was getting really nervous…" "Roxy (neutral smiling_nervous): I haven’t played with everyone else for a while, and I was getting really nervous…"

Replace-ExactText "Roxy: Why were you guys late, though? It’s a little unusual since you usually arrive a lot earlier than everyone else." "Roxy (neutral curious): Why were you guys late, though? It’s a little unusual since you usually arrive a lot earlier than everyone else."

Replace-ExactText "Roxy: Actually, come to think of it you didn’t come on Sunday." "Roxy (neutral thinking): Actually, come to think of it you didn’t come on Sunday."

Replace-ExactText "Roxy: What happened then? Lover’s quarrel?" "Roxy (neutral hehe): What happened then? Lover’s quarrel?"

Replace-ExactText "Roxy: No need to be shy, nothing wrong with high schoolers dating." "Roxy (neutral smiling_eyes_closed): No need to be shy, nothing wrong with high schoolers dating."

Replace-ExactText "Roxy: There are even a lotta middle school couples nowadays, you know." "Roxy (neutral neutral): There are even a lotta middle school couples nowadays, you know."

Replace-ExactText "Roxy: I know, I know…" "Roxy (neutral smiling): I know, I know…"

Replace-ExactText "Roxy: It’s kinda cute when you get flustered, though. I couldn’t help myself." "Roxy (neutral smiling_blushing_eyes): It’s kinda cute when you get flustered, though. I couldn’t help myself."

Replace-ExactText "Roxy: I’m gonna go practice on my own now, so I’ll see you around." "Roxy (neutral smiling_blushing): I’m gonna go practice on my own now, so I’ll see you around."

Replace-ExactText "Roxy: Bye!" "Roxy (waving smiling_blushing): Bye!"

# --- 2) Insert new stand-alone expression-tag paragraphs --------------------
# (inserted after the now-updated anchor lines; bottom-to-top order keeps
#  each Find/anchor text unique and unaffected by earlier insertions)

Insert-TagAfter "Roxy (waving smiling_blushing): Bye!" "Roxy (exit):"

Insert-TagAfter "Pro: That’s not really the point, though…" "Roxy (neutral laughing):"

Insert-TagAfter "Roxy (neutral hehe): What happened then? Lover’s quarrel?" "Roxy (neutral smiling):"

Insert-TagAfter "Pro: Um…" "Roxy (neutral curious):"

Insert-TagAfter "Pro: Ah, some things happened…" "Roxy (neutral skeptical):"

Insert-TagAfter "Pro: Oh, I guess that makes sense." "Roxy (neutral smiling):"

Insert-TagAfter "Pro: Oh, uh, hey." "Roxy (neutral sigh):"

Insert-TagAfter "I start as Roxy suddenly appears to my left." "Roxy (neutral smiling):"

Insert-TagAfter "Roxy (waving smiling): Hey, there." "Roxy (neutral grinning):"

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
